$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "50.118.04"
$ws.Cells.Item(2, 5).Value = "  +4.28%  "

$ws.Cells.Item(3, 4).Value = "2.634.83"
$ws.Cells.Item(3, 5).Value = "  +5.55%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 4).ClearFormats()

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "327.26"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +2.19%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "110.44"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +2.88%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.534"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  +1.82%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -0.11%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.562"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +4.76%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "40.45"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +2.55%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.66"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +2.50%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0822"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +1.32%  "

$ws.Cells.Item(13, 5).Value = "  +0.65%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.29"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +2.48%  "

$ws.Cells.Item(15, 4).Value = "3.041.02"
$ws.Cells.Item(15, 5).Value = "  +5.17%  "

$ws.Cells.Item(16, 4).Value = "2.618.30"
$ws.Cells.Item(16, 5).Value = "  +4.58%  "

$ws.Cells.Item(17, 5).Value = "  +4.96%  "

$ws.Cells.Item(18, 4).Value = "49.970.09"
$ws.Cells.Item(18, 5).Value = "  +4.25%  "

$ws.Cells.Item(19, 5).Value = "  +12.89%  "

$ws.Cells.Item(20, 5).Value = "  +3.09%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.82"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +1.11%  "

$ws.Cells.Item(22, 4).Value = "0.0₃0959"
$ws.Cells.Item(22, 5).Value = "  +2.40%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "72.87"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +1.59%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "279.22"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +1.03%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.60"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +2.17%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.63"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +4.03%  "

$ws.Cells.Item(27, 5).Value = "  -0.01%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "36.91"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +5.67%  "

$ws.Cells.Item(29, 5).Value = "  +5.90%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.95"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +2.42%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.145"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +3.46%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "49.88"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  +1.02%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "19.83"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +1.31%  "

$ws.Cells.Item(34, 5).Value = "  +2.96%  "

$ws.Cells.Item(35, 5).Value = "  -0.26%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0793"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +2.00%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.07"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +6.78%  "

$ws.Cells.Item(38, 5).Value = "  +3.30%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.10"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +7.62%  "

$ws.Cells.Item(40, 5).Value = "  +1.38%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "123.26"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +1.76%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "22.55"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +5.40%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.23"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +0.51%  "

$ws.Cells.Item(44, 5).Value = "  +4.67%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.37"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +7.19%  "

$ws.Cells.Item(46, 4).Value = "2.057.33"
$ws.Cells.Item(46, 5).Value = "  +2.38%  "

$ws.Cells.Item(47, 5).Value = "  +16.23%  "

$ws.Cells.Item(48, 5).Value = "  +9.20%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.05"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +0.82%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.37"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +4.26%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "81.82"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +2.24%  "
